$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Rows 1-4: simple text replacements ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "115"

# Row 5 (0.00003) is unchanged.

# --- Row 6: 0.00008 -> 0.00009 ---
$t.Cell(6,1).Range.Text = "0.00009"

# --- Remove the next three rows (old rows 7,8,9: 0.00005 / 0.00003 / 0.00003) ---
$t.Rows(7).Delete()
$t.Rows(7).Delete()
$t.Rows(7).Delete()

# After the deletions, the row that used to be #10 (0.00004) is now row #7 - unchanged.
# The row that used to be #11 (0.00008) is now row #8 -> becomes 0.00002.
$t.Cell(8,1).Range.Text = "0.00002"
# The row that used to be #12 (0.00015) is now row #9 -> becomes 0.00004.
$t.Cell(9,1).Range.Text = "0.00004"

# --- Insert three new rows before the row that used to be #13 (100.0), now row #10 ---
# (Rows.Add(beforeRow) always inserts immediately above beforeRow, so add them in
# reverse order to end up with 0.00004, 0.00009, 0.00483 reading top to bottom.)
$beforeRow = $t.Rows(10)
$newRow3 = $t.Rows.Add($beforeRow)
$newRow3.Cells(1).Range.Text = "0.00483"
$newRow2 = $t.Rows.Add($beforeRow)
$newRow2.Cells(1).Range.Text = "0.00009"
$newRow1 = $t.Rows.Add($beforeRow)
$newRow1.Cells(1).Range.Text = "0.00004"

# Rows 14-43 (original numbering, still 14-43 after net-zero row count change) are unchanged.

# --- Last three rows: collapse multi-run tab-separated text to a single value ---
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "70"
